$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.405.83"
$ws.Range("E2").Value = "  -5.47%  "

# Row 3
$ws.Range("D3").Value = "'3.473.19"
$ws.Range("E3").Value = "  -6.82%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'558.57"
$ws.Range("E5").Value = "  -8.83%  "

# Row 6
$ws.Range("D6").Value = "'180.52"
$ws.Range("E6").Value = "  -6.62%  "

# Row 7
$ws.Range("D7").Value = "'3.475.93"
$ws.Range("E7").Value = "  -6.63%  "

# Row 8
$ws.Range("E8").Value = "  -6.23%  "

# Row 9
$ws.Range("E9").Value = "  +0.10%  "

# Row 10
$ws.Range("D10").Value = "'0.647"
$ws.Range("E10").Value = "  -11.48%  "

# Row 11
$ws.Range("D11").Value = "'0.141"
$ws.Range("E11").Value = "  -13.15%  "

# Row 12
$ws.Range("D12").Value = "'51.44"
$ws.Range("E12").Value = "  -14.88%  "

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -14.70%  "

# Row 14
$ws.Range("D14").Value = "'9.49"
$ws.Range("E14").Value = "  -10.34%  "

# Row 15
$ws.Range("D15").Value = "'4.030.52"
$ws.Range("E15").Value = "  -6.92%  "

# Row 16
$ws.Range("E16").Value = "  -1.77%  "

# Row 17
$ws.Range("D17").Value = "'3.469.57"
$ws.Range("E17").Value = "  -7.03%  "

# Row 18
$ws.Range("D18").Value = "'65.299.43"
$ws.Range("E18").Value = "  -5.48%  "

# Row 19
$ws.Range("D19").Value = "'17.66"
$ws.Range("E19").Value = "  -9.54%  "

# Row 20
$ws.Range("D20").Value = "'11.68"
$ws.Range("E20").Value = "  -9.74%  "

# Row 21
$ws.Range("E21").Value = "  -11.02%  "

# Row 22
$ws.Range("D22").Value = "'377.53"
$ws.Range("E22").Value = "  -8.67%  "

# Row 23
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  -11.31%  "

# Row 24
$ws.Range("D24").Value = "'83.33"
$ws.Range("E24").Value = "  -7.65%  "

# Row 25
$ws.Range("D25").Value = "'10.58"
$ws.Range("E25").Value = "  -6.35%  "

# Row 26
$ws.Range("E26").Value = "  -8.95%  "

# Row 27
$ws.Range("E27").Value = "  -1.18%  "

# Row 28
$ws.Range("E28").Value = "  -8.06%  "

# Row 29
$ws.Range("D29").Value = "'3.46"
$ws.Range("E29").Value = "  -9.10%  "

# Row 30
$ws.Range("D30").Value = "'8.62"
$ws.Range("E30").Value = "  -11.90%  "

# Row 31
$ws.Range("D31").Value = "'30.37"
$ws.Range("E31").Value = "  -8.06%  "

# Row 32
$ws.Range("D32").Value = "'7.19"
$ws.Range("E32").Value = "  -7.61%  "

# Row 33
$ws.Range("D33").Value = "'11.90"
$ws.Range("E33").Value = "  -7.12%  "

# Row 34
$ws.Range("D34").Value = "'604.05"
$ws.Range("E34").Value = "  -5.79%  "

# Row 35
$ws.Range("D35").Value = "'62.11"
$ws.Range("E35").Value = "  -7.78%  "

# Row 36
$ws.Range("E36").Value = "  -10.59%  "

# Row 37
$ws.Range("D37").Value = "'40.50"
$ws.Range("E37").Value = "  -11.88%  "

# Row 38
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("D39").Value = "'0.390"
$ws.Range("E39").Value = "  -6.40%  "

# Row 40
$ws.Range("E40").Value = "  -0.16%  "

# Row 41
$ws.Range("D41").Value = "'0.0₃0709"
$ws.Range("E41").Value = "  -15.04%  "

# Row 42
$ws.Range("E42").Value = "  -9.11%  "

# Row 43
$ws.Range("D43").Value = "'2.947.79"
$ws.Range("E43").Value = "  +1.43%  "

# Row 44
$ws.Range("E44").Value = "  -11.54%  "

# Row 45
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  -8.18%  "

# Row 46
$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  +0.34%  "

# Row 47
$ws.Range("E47").Value = "  -12.38%  "

# Row 48
$ws.Range("D48").Value = "'0.127"
$ws.Range("E48").Value = "  -9.74%  "

# Row 49
$ws.Range("D49").Value = "'135.96"
$ws.Range("E49").Value = "  -4.78%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'8.13"
$ws.Range("E50").Value = "  -12.22%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.66"
$ws.Range("E51").Value = "  -4.57%  "

